# Weekly data refresh: a new record for the current week is inserted at the
# top of the data table (row 47, just after the still-unique first block of
# rows 2-46), pushing every subsequent record down by one row.
#
# Sheet columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
# F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
# K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
# N Unidad de comercialización, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificación.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 47; Excel shifts the old row 47..152 down to
# 48..153 and extends the used range / dimension to R153 automatically.
$ws.Rows(47).Insert()

$newRow = 47
$ws.Cells.Item($newRow, 1).Value  = 11
$ws.Cells.Item($newRow, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item($newRow, 3).Value  = "Bíobío"
$ws.Cells.Item($newRow, 4).Value  = 45259
$ws.Cells.Item($newRow, 5).Value  = 8
$ws.Cells.Item($newRow, 6).Value  = 100112037
$ws.Cells.Item($newRow, 7).Value  = "Cebollín"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 50
$ws.Cells.Item($newRow, 11).Value = 3500
$ws.Cells.Item($newRow, 12).Value = 3500
$ws.Cells.Item($newRow, 13).Value = 3500
$ws.Cells.Item($newRow, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item($newRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value = 97
$ws.Cells.Item($newRow, 17).Value = 36
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"

# Match the date formatting already used by the rest of column D.
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
